# Form the consolidated report: set the "Absent" (column H) values for
# the first two date rows (rows 3 and 4) and the corresponding rows 9 and 10
# further down the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 0
